{"js": "// Replace each three-digit \u00f7 one-digit practice answer in the document's\n// table with its new value, per the author's regenerated problem set.\n// Every \"old\" string below is unique in the document (single occurrence),\n// so a scoped search + full-range text replace is safe and unambiguous.\nconst replacements = [\n  [\"463\u00f78=57, 7\", \"135\u00f73=45, 0\"],\n  [\"759\u00f73=253, 0\", \"284\u00f77=40, 4\"],\n  [\"886\u00f78=110, 6\", \"515\u00f78=64, 3\"],\n  [\"120\u00f75=24, 0\", \"434\u00f75=86, 4\"],\n  [\"606\u00f77=86, 4\", \"307\u00f74=76, 3\"],\n  [\"965\u00f78=120, 5\", \"869\u00f79=96, 5\"],\n  [\"977\u00f79=108, 5\", \"930\u00f76=155, 0\"],\n  [\"964\u00f75=192, 4\", \"950\u00f79=105, 5\"],\n  [\"449\u00f75=89, 4\", \"362\u00f75=72, 2\"],\n  [\"378\u00f73=126, 0\", \"779\u00f79=86, 5\"],\n  [\"311\u00f78=38, 7\", \"158\u00f76=26, 2\"],\n  [\"729\u00f77=104, 1\", \"249\u00f78=31, 1\"],\n  [\"108\u00f73=36, 0\", \"908\u00f75=181, 3\"],\n  [\"746\u00f78=93, 2\", \"440\u00f73=146, 2\"],\n  [\"256\u00f76=42, 4\", \"410\u00f78=51, 2\"],\n  [\"139\u00f77=19, 6\", \"436\u00f74=109, 0\"],\n  [\"321\u00f74=80, 1\", \"167\u00f74=41, 3\"],\n  [\"297\u00f79=33, 0\", \"400\u00f79=44, 4\"],\n  [\"995\u00f78=124, 3\", \"805\u00f75=161, 0\"],\n  [\"554\u00f79=61, 5\", \"672\u00f72=336, 0\"],\n  [\"928\u00f78=116, 0\", \"190\u00f78=23, 6\"],\n  [\"743\u00f79=82, 5\", \"985\u00f79=109, 4\"],\n  [\"330\u00f72=165, 0\", \"198\u00f79=22, 0\"],\n  [\"193\u00f74=48, 1\", \"315\u00f72=157, 1\"],\n  [\"833\u00f73=277, 2\", \"279\u00f78=34, 7\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each three-digit \u00f7 one-digit practice answer in the document's\n# table with its new value, per the author's regenerated problem set.\n# Every \"old\" string is unique in the document (single occurrence), so a\n# Find/Replace over the whole story content is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"463\u00f78=57, 7\", \"135\u00f73=45, 0\"),\n  @(\"759\u00f73=253, 0\", \"284\u00f77=40, 4\"),\n  @(\"886\u00f78=110, 6\", \"515\u00f78=64, 3\"),\n  @(\"120\u00f75=24, 0\", \"434\u00f75=86, 4\"),\n  @(\"606\u00f77=86, 4\", \"307\u00f74=76, 3\"),\n  @(\"965\u00f78=120, 5\", \"869\u00f79=96, 5\"),\n  @(\"977\u00f79=108, 5\", \"930\u00f76=155, 0\"),\n  @(\"964\u00f75=192, 4\", \"950\u00f79=105, 5\"),\n  @(\"449\u00f75=89, 4\", \"362\u00f75=72, 2\"),\n  @(\"378\u00f73=126, 0\", \"779\u00f79=86, 5\"),\n  @(\"311\u00f78=38, 7\", \"158\u00f76=26, 2\"),\n  @(\"729\u00f77=104, 1\", \"249\u00f78=31, 1\"),\n  @(\"108\u00f73=36, 0\", \"908\u00f75=181, 3\"),\n  @(\"746\u00f78=93, 2\", \"440\u00f73=146, 2\"),\n  @(\"256\u00f76=42, 4\", \"410\u00f78=51, 2\"),\n  @(\"139\u00f77=19, 6\", \"436\u00f74=109, 0\"),\n  @(\"321\u00f74=80, 1\", \"167\u00f74=41, 3\"),\n  @(\"297\u00f79=33, 0\", \"400\u00f79=44, 4\"),\n  @(\"995\u00f78=124, 3\", \"805\u00f75=161, 0\"),\n  @(\"554\u00f79=61, 5\", \"672\u00f72=336, 0\"),\n  @(\"928\u00f78=116, 0\", \"190\u00f78=23, 6\"),\n  @(\"743\u00f79=82, 5\", \"985\u00f79=109, 4\"),\n  @(\"330\u00f72=165, 0\", \"198\u00f79=22, 0\"),\n  @(\"193\u00f74=48, 1\", \"315\u00f72=157, 1\"),\n  @(\"833\u00f73=277, 2\", \"279\u00f78=34, 7\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n\n  $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
